# Apply "Penalty Reward System" edits (unfinished, per commit message):
# 1. Weekly Quantity sheet: remove the row for week 45368.99999999999 (old row 4),
#    shifting all subsequent weekly rows up by one (dimension becomes A1:B24).
# 2. Monthly Trend sheet: update the requested quantity for the
#    45382.99999999999 month (row 3) from 240 to 80.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows(4).Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B3").Value = 80
